$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was the "Perfection" / Provincia de Huasco / 44335 row) -> becomes the
# "Sin especificar" / Región del Maule / 44162 row
$ws.Range("D2").Value = 44162
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("N2").Value = "$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 700

# Row 3 (was 44454 / Provincia de Limarí) -> becomes 44342 row
$ws.Range("D3").Value = 44342
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 30000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 31000
$ws.Range("P3").Value = 1240

# Row 4 (was 44342 / Provincia de Limarí) -> becomes 44399 / Provincia de Huasco row
$ws.Range("D4").Value = 44399
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 39000
$ws.Range("L4").Value = 40000
$ws.Range("M4").Value = 39600
$ws.Range("O4").Value = "Provincia de Huasco"
$ws.Range("P4").Value = 1584

# Row 5 (was 44399 / Provincia de Huasco) -> becomes 44454 / Provincia de Limarí row
$ws.Range("D5").Value = 44454
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 36000
$ws.Range("L5").Value = 38000
$ws.Range("M5").Value = 37000
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 1480

# Row 6 (was 44328 / Provincia de Huasco) -> becomes 44335 row
$ws.Range("D6").Value = 44335
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 32000
$ws.Range("M6").Value = 31000
$ws.Range("P6").Value = 1240

# Row 7 (was 44162 / Sin especificar / Región del Maule) -> becomes 44328 / Perfection / Provincia de Huasco row
$ws.Range("D7").Value = 44328
$ws.Range("H7").Value = "Perfection"
$ws.Range("K7").Value = 33000
$ws.Range("L7").Value = 34000
$ws.Range("M7").Value = 33500
$ws.Range("N7").Value = "$/malla 25 kilos"
$ws.Range("O7").Value = "Provincia de Huasco"
$ws.Range("P7").Value = 1340
